# Clarify learner roles in masking workflow
$p = $ppt.ActivePresentation

# Slide 5 - AC測定プロセス / 操作フロー
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(5).TextFrame.TextRange.Text = "2. 受講者が 5dB上昇法（聴取不可レベルから 5dB ずつ上げる）で応答を確認"

# Slide 6 - BC測定 & ABG管理 / 測定時の確認ポイント
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(8).TextFrame.TextRange.Text = "- 受講者が上昇法で応答を記録し、仮閾値と AC の差から ABG を算出"

# Slide 7 - マスキング要否判定ロジック
$s7 = $p.Slides.Item(7)
$s7.Shapes.Item(3).TextFrame.TextRange.Text = "条件評価（受講者が判断）"
$s7.Shapes.Item(4).TextFrame.TextRange.Text = "- ABG が疾患別最小値以上か"
$s7.Shapes.Item(5).TextFrame.TextRange.Text = "- 左右差が IA を超える見込みか"
$s7.Shapes.Item(6).TextFrame.TextRange.Text = "- テスト耳閾値がマスキングノイズで覆われないか"
$s7.Shapes.Item(7).TextFrame.TextRange.Text = "システムのサポート（AI）"
$s7.Shapes.Item(9).TextFrame.TextRange.Text = "- 必要なマスキング操作をガイダンス"
$s7.Shapes.Item(10).TextFrame.TextRange.Text = "- 判定ミス時に根拠をフィードバック"

# Slide 8 - マスキング量算出
$s8 = $p.Slides.Item(8)
$s8.Shapes.Item(3).TextFrame.TextRange.Text = "推奨計算（受講者が手計算）"
$s8.Shapes.Item(6).TextFrame.TextRange.Text = "- セーフティマージンは 10dB を基準に調整"
$s8.Shapes.Item(7).TextFrame.TextRange.Text = "アプリでの表示（AI サポート）"
$s8.Shapes.Item(9).TextFrame.TextRange.Text = "- 入力値と推奨値の乖離をリアルタイム判定"

# Slide 11 - 全体フローチャート (Mermaid source lines)
$s11 = $p.Slides.Item(11)
$s11.Shapes.Item(4).TextFrame.TextRange.Text = "  subgraph AI[AIの役割]"
$s11.Shapes.Item(5).TextFrame.TextRange.Text = "    Start[症例生成] --> ACGen[AC/BC初期値と補助検査生成]"
$s11.Shapes.Item(6).TextFrame.TextRange.Text = "    ACGen --> Guidance[判定ガイド・推奨値提示]"
$s11.Shapes.Item(7).TextFrame.TextRange.Text = "    Feedback[結果照合・フィードバック] --> Start"
$s11.Shapes.Item(8).TextFrame.TextRange.Text = "  end"
$s11.Shapes.Item(9).TextFrame.TextRange.Text = "  subgraph Learner[受講者の役割]"
$s11.Shapes.Item(10).TextFrame.TextRange.Text = "    ACMeasure[AC測定 5dB上昇法] --> BCMeasure[BC測定とABG管理]"
$s11.Shapes.Item(11).TextFrame.TextRange.Text = "    BCMeasure --> MaskCheck[マスキング要否判定]"
$s11.Shapes.Item(12).TextFrame.TextRange.Text = "    MaskCheck -->|必要| MaskCalc[マスキング量算出]"
